$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Objetivos: / Objectives text cells now hold the instructor reference
# that used to live further down the sheet.
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# Row 13 ("Programa resumido:") B/C must contain the same text as B8/C8
# ("01/01/2023"). Copy from B8/C8 instead of typing the literal string so
# that Excel does not reinterpret it as a date serial value and does not
# introduce a new number-format style.
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4104) | Out-Null
$ws.Range("C8").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

# Row 15 ("Programa:") B/C now hold the instructor reference as well.
$ws.Range("B15").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C15").Value = "519033 - Carlos Yujiro Shigue"

# Row 18 ("Método:") B/C now hold the second instructor reference.
$ws.Range("B18").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C18").Value = "7290967 - Emerson Gonçalves de Melo"
